$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C17").Value = -7.35543964232488
$ws.Range("D17").Value = 4.3
$ws.Range("E17").Value = -15.9
